# Apply the two functional changes captured by the commit:
#
# 1. The table on slide 6 (the "SOURCES OF FINANCE" table) is switched to a
#    different built-in PowerPoint table style
#    ({0FE91767-8CCD-4ABC-8105-AB7A362103CE}).  Table styles are applied
#    through Table.ApplyStyle(styleId) - assigning Table.Style directly is
#    rejected by the object model.
#
# 2. The presentation's applied theme ("Integral") colours are swapped back
#    to the default Office theme colours.  The theme's colour slots are
#    exposed on Master.Theme.ThemeColorScheme as a fixed, ordered set of 12
#    entries (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink); re-pointing
#    each entry at the stock "Office" RGB values reproduces the colour
#    swap that the commit performs on the theme part backing the slide
#    master.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $tableShape = $shp
        break
    }
}
$table = $tableShape.Table
$table.ApplyStyle("{0FE91767-8CCD-4ABC-8105-AB7A362103CE}")

# --- 2. Theme colours -------------------------------------------------------
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
$colorScheme.Item(1).RGB  = 0        # dk1      000000
$colorScheme.Item(2).RGB  = 16777215 # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 6968388  # dk2      44546A
$colorScheme.Item(4).RGB  = 15132391 # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 13998939 # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 3243501  # accent2  ED7D31
$colorScheme.Item(7).RGB  = 10855845 # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 49407    # accent4  FFC000
$colorScheme.Item(9).RGB  = 12874308 # accent5  4472C4
$colorScheme.Item(10).RGB = 4697456  # accent6  70AD47
$colorScheme.Item(11).RGB = 12673797 # hlink    0563C1
$colorScheme.Item(12).RGB = 7491477  # folHlink 954F72
